$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020408781130923
$ws.Range("D2").Value = 1.02878434037026
$ws.Range("E2").Value = 1.030707764241218
$ws.Range("F2").Value = 1.040100143588177
$ws.Range("I2").Value = 1.030696972397081
$ws.Range("J2").Value = 1.0256060573022
$ws.Range("K2").Value = 1.031600150515
$ws.Range("L2").Value = 1.033517994970504
$ws.Range("M2").Value = 1.042883451044268
$ws.Range("N2").Value = 1.012560791950009
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021393349247925
$ws.Range("D3").Value = 1.02914128211441
$ws.Range("E3").Value = 1.031585593865843
$ws.Range("F3").Value = 1.041071780442473
$ws.Range("I3").Value = 1.030723332733326
$ws.Range("J3").Value = 1.026227617167126
$ws.Range("K3").Value = 1.03176631719081
$ws.Range("L3").Value = 1.034204049394636
$ws.Range("M3").Value = 1.043665017484528
$ws.Range("N3").Value = 1.012771732068521
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022031139767109
$ws.Range("D4").Value = 1.029372071347479
$ws.Range("E4").Value = 1.032154518113247
$ws.Range("F4").Value = 1.041701271679765
$ws.Range("I4").Value = 1.030738892160625
$ws.Range("J4").Value = 1.026629947757843
$ws.Range("K4").Value = 1.031872913766997
$ws.Range("L4").Value = 1.03464824669775
$ws.Range("M4").Value = 1.044170903434631
$ws.Range("N4").Value = 1.012908136370038
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022299436029145
$ws.Range("D5").Value = 1.029469050915985
$ws.Range("E5").Value = 1.032393910499154
$ws.Range("F5").Value = 1.041966094538081
$ws.Range("I5").Value = 1.030745074042366
$ws.Range("J5").Value = 1.026799120254467
$ws.Range("K5").Value = 1.031917504295126
$ws.Range("L5").Value = 1.034835051958904
$ws.Range("M5").Value = 1.044383615039873
$ws.Range("N5").Value = 1.012965459295459
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022344494051537
$ws.Range("D6").Value = 1.029485331557471
$ws.Range("E6").Value = 1.032434118203549
$ws.Range("F6").Value = 1.042010570256145
$ws.Range("I6").Value = 1.030746090914481
$ws.Range("J6").Value = 1.02682752695777
$ws.Range("K6").Value = 1.031924978143245
$ws.Range("L6").Value = 1.03486642115208
$ws.Range("M6").Value = 1.044419332398107
$ws.Range("N6").Value = 1.012975082798957
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022034724091119
$ws.Range("D7").Value = 1.029373367369667
$ws.Range("E7").Value = 1.032157716036278
$ws.Range("F7").Value = 1.041704809531519
$ws.Range("I7").Value = 1.030738976176053
$ws.Range("J7").Value = 1.02663220812141
$ws.Range("K7").Value = 1.031873510464511
$ws.Range("L7").Value = 1.034650742545679
$ws.Range("M7").Value = 1.044173745553044
$ws.Range("N7").Value = 1.012908902406481
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02074137332579
$ws.Range("D8").Value = 1.028905005672677
$ws.Range("E8").Value = 1.031004241791601
$ws.Range("F8").Value = 1.040428351085138
$ws.Range("I8").Value = 1.030706190575577
$ws.Range("J8").Value = 1.025816086970659
$ws.Range("K8").Value = 1.031656497894039
$ws.Range("L8").Value = 1.033749792935448
$ws.Range("M8").Value = 1.04314755118848
$ws.Range("N8").Value = 1.012632098166917
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018467783051723
$ws.Range("D9").Value = 1.028078437634381
$ws.Range("E9").Value = 1.028978692363243
$ws.Range("F9").Value = 1.038185079370054
$ws.Range("I9").Value = 1.030636984707148
$ws.Range("J9").Value = 1.024379081062558
$ws.Range("K9").Value = 1.031267074970936
$ws.Range("L9").Value = 1.032164345964668
$ws.Range("M9").Value = 1.041340545575237
$ws.Range("N9").Value = 1.012143675226395
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016955762872941
$ws.Range("D10").Value = 1.027526673793902
$ws.Range("E10").Value = 1.027633114536495
$ws.Range("F10").Value = 1.036693684458129
$ws.Range("I10").Value = 1.030583208002578
$ws.Range("J10").Value = 1.023421863183637
$ws.Range("K10").Value = 1.031002822270585
$ws.Range("L10").Value = 1.031108875798346
$ws.Range("M10").Value = 1.040136806409196
$ws.Range("N10").Value = 1.011817638421183
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016301927770907
$ws.Range("D11").Value = 1.027287610094918
$ws.Range("E11").Value = 1.027051614499201
$ws.Range("F11").Value = 1.036048887236079
$ws.Range("I11").Value = 1.030558120182267
$ws.Range("J11").Value = 1.023007572705944
$ws.Range("K11").Value = 1.030887315303973
$ws.Range("L11").Value = 1.030652211324026
$ws.Range("M11").Value = 1.03961580956776
$ws.Range("N11").Value = 1.011676365300691
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016059196967308
$ws.Range("D12").Value = 1.027198791374191
$ws.Range("E12").Value = 1.026835792493957
$ws.Range("F12").Value = 1.035809530330271
$ws.Range("I12").Value = 1.030548531491896
$ws.Range("J12").Value = 1.022853716187377
$ws.Range("K12").Value = 1.030844249440951
$ws.Range("L12").Value = 1.030482640929205
$ws.Range("M12").Value = 1.039422323864938
$ws.Range("N12").Value = 1.01162387593721
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016111257536946
$ws.Range("D13").Value = 1.027217844153525
$ws.Range("E13").Value = 1.026882079200083
$ws.Range("F13").Value = 1.035860866413929
$ws.Range("I13").Value = 1.030550600501898
$ws.Range("J13").Value = 1.02288671759089
$ws.Range("K13").Value = 1.030853494497589
$ws.Range("L13").Value = 1.030519011845939
$ws.Range("M13").Value = 1.03946382557016
$ws.Range("N13").Value = 1.011635135722592
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016281860857633
$ws.Range("D14").Value = 1.027280268703587
$ws.Range("E14").Value = 1.027033771049288
$ws.Range("F14").Value = 1.036029098862666
$ws.Range("I14").Value = 1.030557333080575
$ws.Range("J14").Value = 1.022994854274154
$ws.Range("K14").Value = 1.030883758751432
$ws.Range("L14").Value = 1.030638193452703
$ws.Range("M14").Value = 1.039599815243887
$ws.Range("N14").Value = 1.011672026801971
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016386992789059
$ws.Range("D15").Value = 1.02731872796747
$ws.Range("E15").Value = 1.027127256361768
$ws.Range("F15").Value = 1.036132772273709
$ws.Range("I15").Value = 1.030561445496355
$ws.Range("J15").Value = 1.023061484771576
$ws.Range("K15").Value = 1.030902384207306
$ws.Range("L15").Value = 1.030711632504841
$ws.Range("M15").Value = 1.039683607737742
$ws.Range("N15").Value = 1.011694754737308
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016999174331396
$ws.Range("D16").Value = 1.027542536766567
$ws.Range("E16").Value = 1.027671730978089
$ws.Range("F16").Value = 1.036736498420884
$ws.Range("I16").Value = 1.030584835097339
$ws.Range("J16").Value = 1.02344936236486
$ws.Range("K16").Value = 1.031010465372445
$ws.Range("L16").Value = 1.0311391907977
$ws.Range("M16").Value = 1.040171388200393
$ws.Range("N16").Value = 1.011827012237716
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017383415496877
$ws.Range("D17").Value = 1.027682888472178
$ws.Range("E17").Value = 1.028013572602913
$ws.Range("F17").Value = 1.037115464919664
$ws.Range("I17").Value = 1.030599024851759
$ws.Range("J17").Value = 1.023692719401159
$ws.Range("K17").Value = 1.0310779726553
$ws.Range("L17").Value = 1.031407484106648
$ws.Range("M17").Value = 1.04047742242415
$ws.Range("N17").Value = 1.011909948136339
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017607621702939
$ws.Range("D18").Value = 1.027764738985853
$ws.Range("E18").Value = 1.028213073347276
$ws.Range("F18").Value = 1.037336604776038
$ws.Range("I18").Value = 1.030607127560651
$ws.Range("J18").Value = 1.02383468384298
$ws.Range("K18").Value = 1.031117243846707
$ws.Range("L18").Value = 1.031564009833938
$ws.Range("M18").Value = 1.040655949139614
$ws.Range("N18").Value = 1.0119583138577
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017684084599193
$ws.Range("D19").Value = 1.02779264541077
$ws.Range("E19").Value = 1.028281116621651
$ws.Range("F19").Value = 1.037412023846064
$ws.Range("I19").Value = 1.030609860848365
$ws.Range("J19").Value = 1.02388309314898
$ws.Range("K19").Value = 1.03113061651259
$ws.Range("L19").Value = 1.031617386934498
$ws.Range("M19").Value = 1.040716825855034
$ws.Range("N19").Value = 1.011974803716529
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017342181277592
$ws.Range("D20").Value = 1.027667831521965
$ws.Range("E20").Value = 1.027976884813917
$ws.Range("F20").Value = 1.037074795554083
$ws.Range("I20").Value = 1.030597520410542
$ws.Range("J20").Value = 1.023666607586091
$ws.Range("K20").Value = 1.031070740578591
$ws.Range("L20").Value = 1.031378695169946
$ws.Range("M20").Value = 1.040444585550877
$ws.Range("N20").Value = 1.011901050867587
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016231618765494
$ws.Range("D21").Value = 1.027261886759268
$ws.Range("E21").Value = 1.026989096795862
$ws.Range("F21").Value = 1.035979554462718
$ws.Range("I21").Value = 1.030555357949713
$ws.Range("J21").Value = 1.022963009898734
$ws.Range("K21").Value = 1.030874851128149
$ws.Range("L21").Value = 1.03060309591129
$ws.Range("M21").Value = 1.039559768679155
$ws.Range("N21").Value = 1.011661163688827
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015534131856222
$ws.Range("D22").Value = 1.027006539571792
$ws.Range("E22").Value = 1.026369036744021
$ws.Range("F22").Value = 1.035291798461181
$ws.Range("I22").Value = 1.030527287241729
$ws.Range("J22").Value = 1.02252080061358
$ws.Range("K22").Value = 1.030750754341557
$ws.Range("L22").Value = 1.030115765239667
$ws.Range("M22").Value = 1.03900365643402
$ws.Range("N22").Value = 1.011510254900703
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015903809909862
$ws.Range("D23").Value = 1.027141913996602
$ws.Range("E23").Value = 1.026697647005945
$ws.Range("F23").Value = 1.035656308464185
$ws.Range("I23").Value = 1.03054231580187
$ws.Range("J23").Value = 1.022755207729864
$ws.Range("K23").Value = 1.030816628399333
$ws.Range("L23").Value = 1.03037407788376
$ws.Range("M23").Value = 1.039298442031205
$ws.Range("N23").Value = 1.011590262181813
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017360812984233
$ws.Range("D24").Value = 1.027674635154001
$ws.Range("E24").Value = 1.027993462105811
$ws.Range("F24").Value = 1.037093171996013
$ws.Range("I24").Value = 1.030598200740298
$ws.Range("J24").Value = 1.023678406334533
$ws.Range("K24").Value = 1.031074008766332
$ws.Range("L24").Value = 1.031391703543831
$ws.Range("M24").Value = 1.04045942305216
$ws.Range("N24").Value = 1.011905071189384
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019054909995131
$ws.Range("D25").Value = 1.028292261044263
$ws.Range("E25").Value = 1.029501506495714
$ws.Range("F25").Value = 1.038764298498148
$ws.Range("I25").Value = 1.030656226074656
$ws.Range("J25").Value = 1.024750446210349
$ws.Range("K25").Value = 1.031267074970936
$ws.Range("L25").Value = 1.032573962613719
$ws.Range("M25").Value = 1.041807540882721
$ws.Range("N25").Value = 1.012270020113531
